# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# 1) El trabajador EUSEBIO VILLAR GUERRERO (filas 23-26) se retira del
#    estado de cuenta.
# 2) Los periodos de mora de CRISTIAM MORALES MUÑOZ (filas 16-22) se
#    reordenan de forma ascendente (2310 -> 2404).
# 3) Los totales (Valor Mora, Cant. Trabajadores, Cant. Periodos) se
#    actualizan para reflejar la nueva base de datos.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Quitar las filas del trabajador que ya no aplica ---------------
$ws.Range("B23:J26").EntireRow.Delete()

# --- 2) Reordenar el detalle restante por Periodo Mora (ascendente) ----
$dataRange = $ws.Range("B16:J22")
$keyRange = $ws.Range("E16:E22")
$dataRange.Sort($keyRange, 1)

# --- 3) Actualizar los totales del encabezado ---------------------------
$ws.Range("E11").Value = 313974
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 7
